$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Sequence values for the last two rows (D3 and D4)
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 1

# Update the active selection to D4
$ws.Range("D4").Select()
